# Actualización automática 2025-07-22 17:20:08
#
# The sale amount for "240X80 PORCELANATO" for advisor "LOZANO MOLINA TITO" /
# client "PAREDES ORTIZ MARIA INES" for the month of "julio" (July)
# increased from 1900.8 to 6514.56 (an increase of 4613.76).
# This value flows through three sheets that need to be kept consistent:
#   - "VENTAS POR GRUPO"    : D18 (240X80 PORCELANATO column for this row)
#   - "VENTA MENSUAL"       : F18 (julio column for this row) and F29 (total)
#   - "CUMPLIMIENTO MENSUAL": D3 (VENTA for 240X80 PORCELANATO) and D19 (TOTAL),
#                             plus the dependent E (POR CUMPLIR = C - D) and
#                             F (CUMPLIMIENTO = D / C) columns on those rows.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D18").Value = 6514.56

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F18").Value = 8691.84
$wsMensual.Range("F29").Value = 10744.08

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO totals for LOZANO MOLINA TITO
$wsCumpl.Range("D3").Value = 6514.56
$wsCumpl.Range("E3").Value = -3394.4455
$wsCumpl.Range("F3").Value = 2.08792337588893

# Row 19: TOTAL row
$wsCumpl.Range("D19").Value = 10978.55
$wsCumpl.Range("E19").Value = 16202.76093005039
$wsCumpl.Range("F19").Value = 0.4039006811795317
